$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 127, shifting existing rows 127:240 down to 128:241
$ws.Rows("127:127").Insert()

# Populate the newly inserted row 127 with the new record
$ws.Cells.Item(127,1).Value2  = 5
$ws.Cells.Item(127,2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(127,3).Value2  = "Maule"
$ws.Cells.Item(127,4).Value2  = 44669
$ws.Cells.Item(127,5).Value2  = 7
$ws.Cells.Item(127,6).Value2  = 100112021
$ws.Cells.Item(127,7).Value2  = "Ají"
$ws.Cells.Item(127,8).Value2  = "Cristal"
$ws.Cells.Item(127,9).Value2  = "Primera"
$ws.Cells.Item(127,10).Value2 = 100
$ws.Cells.Item(127,11).Value2 = 15000
$ws.Cells.Item(127,12).Value2 = 15000
$ws.Cells.Item(127,13).Value2 = 15000
$ws.Cells.Item(127,14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(127,15).Value2 = "Región del Maule"
$ws.Cells.Item(127,16).Value2 = 600
$ws.Cells.Item(127,17).Value2 = 25
$ws.Cells.Item(127,18).Value2 = "Hortaliza"
